$d = $word.ActiveDocument

# Replace the "Original File Name" placeholder value with the actual
# supplied original file name.
$d.Content.Find.Execute(
    "Original File Name", $true, $false, $false, $false, $false,
    $true, 1, $false, "Diplome-dhematologie-Original", 2) | Out-Null

# Correct the certified date from March 13, 2023 to March 15, 2023.
$d.Content.Find.Execute(
    "March 13, 2023", $true, $false, $false, $false, $false,
    $true, 1, $false, "March 15, 2023", 2) | Out-Null
